# Update the "Förändrad" (Changed) date column (C) for rows 2-28
# from 2024-06-24 (serial 45467) to 2024-06-25 (serial 45468).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45467) {
        $cell.Value = 45468
    }
}
